$wb = $excel.ActiveWorkbook

# --- "studies" sheet: add a new "PMID" header column (H1) ---
$studies = $wb.Worksheets.Item("studies")
$studies.Activate()
$studies.Range("H1").Value = "PMID"
$studies.Range("H2").Select()

# --- "counts" sheet: add a new "notes" header column (F1) ---
$counts = $wb.Worksheets.Item("counts")
$counts.Activate()
$counts.Range("F1").Value = "notes"
$counts.Range("F2").Select()

# "counts" ends up as the active sheet/tab after the edit
$counts.Activate()
